# Add two new product sheets (ID_867f88d, ID_49d69c3) to the workbook,
# matching the structure/style used by the existing product sheets
# (header row: price | date | name | description | id).

$wb = $excel.ActiveWorkbook

$headers = @("price", "date", "name", "description", "id")

$newSheets = @(
    @{
        Name = "ID_867f88d"
        Price = 15.99
        Date = "17/02/2025"
        ProductName = "Mens Casual Slim Fit"
        Description = "The color could be slightly different between on the screen and in practice. / Please note that body builds vary by person, therefore, detailed size information should be reviewed below on the product description."
        Id = "ID_867f88d"
    },
    @{
        Name = "ID_49d69c3"
        Price = 9.99
        Date = "17/02/2025"
        ProductName = "White Gold Plated Princess"
        Description = "Classic Created Wedding Engagement Solitaire Diamond Promise Ring for Her. Gifts to spoil your love more for Engagement, Wedding, Anniversary, Valentine's Day..."
        Id = "ID_49d69c3"
    }
)

foreach ($sheetInfo in $newSheets) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $sheetInfo.Name

    # Header row
    $headerRange = $ws.Range("A1:E1")
    $ws.Range("A1").Value = $headers[0]
    $ws.Range("B1").Value = $headers[1]
    $ws.Range("C1").Value = $headers[2]
    $ws.Range("D1").Value = $headers[3]
    $ws.Range("E1").Value = $headers[4]

    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108  # xlCenter
    $headerRange.VerticalAlignment = -4160    # xlTop
    $headerRange.Borders.LineStyle = 1
    $headerRange.Borders.Weight = 2

    # Data row
    $ws.Range("A2").Value = $sheetInfo.Price
    $ws.Range("B2").Value = $sheetInfo.Date
    $ws.Range("C2").Value = $sheetInfo.ProductName
    $ws.Range("D2").Value = $sheetInfo.Description
    $ws.Range("E2").Value = $sheetInfo.Id
}

Write-Output "Added sheets: ID_867f88d, ID_49d69c3"
